$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.890.76'
$ws.Range("E2").Value = '  -0.66%  '
$ws.Range("D3").Value = '1.861.01'
$ws.Range("E3").Value = '  -0.27%  '
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = "'304.68"
$ws.Range("E5").Value = '  -0.59%  '
$ws.Range("D6").Value = "'0.9997"
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("D7").Value = "'0.5045"
$ws.Range("E7").Value = '  -1.03%  '
$ws.Range("D8").Value = "'0.3637"
$ws.Range("E8").Value = '  -2.83%  '
$ws.Range("D9").Value = "'0.07182"
$ws.Range("E9").Value = '  +0.72%  '
$ws.Range("D10").Value = "'0.8950"
$ws.Range("E10").Value = '  +0.87%  '
$ws.Range("D11").Value = "'20.76"
$ws.Range("E11").Value = '  +0.71%  '
$ws.Range("D12").Value = '1.863.52'
$ws.Range("E12").Value = '  -0.13%  '
$ws.Range("D13").Value = "'0.07493"
$ws.Range("E13").Value = '  -0.32%  '
$ws.Range("D14").Value = "'92.63"
$ws.Range("E14").Value = '  +4.06%  '
$ws.Range("D15").Value = "'5.228"
$ws.Range("E15").Value = '  -1.52%  '
$ws.Range("D16").Value = "'0.9996"
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("D17").Value = "'0.000008489"
$ws.Range("E17").Value = '  +0.20%  '
$ws.Range("D18").Value = "'14.19"
$ws.Range("E18").Value = '  +0.53%  '
$ws.Range("D19").Value = "'0.9993"
$ws.Range("E19").Value = '  +0.01%  '
$ws.Range("D20").Value = '26.927.17'
$ws.Range("E20").Value = '  -0.74%  '
$ws.Range("D21").Value = "'5.036"
$ws.Range("E21").Value = '  -0.15%  '
$ws.Range("D22").Value = '2.105.33'
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("D23").Value = "'10.36"
$ws.Range("E23").Value = '  -1.70%  '
$ws.Range("D24").Value = "'6.403"
$ws.Range("E24").Value = '  -1.08%  '
$ws.Range("D25").Value = "'147.49"
$ws.Range("E25").Value = '  -1.38%  '
$ws.Range("D26").Value = "'1.781"
$ws.Range("E26").Value = '  -3.56%  '
$ws.Range("D27").Value = "'17.89"
$ws.Range("E27").Value = '  -0.19%  '
$ws.Range("D28").Value = "'2.064"
$ws.Range("E28").Value = '  -1.68%  '
$ws.Range("D29").Value = "'113.09"
$ws.Range("E29").Value = '  +0.30%  '
$ws.Range("D30").Value = "'4.679"
$ws.Range("E30").Value = '  -1.01%  '
$ws.Range("D31").Value = "'4.680"
$ws.Range("E31").Value = '  -0.07%  '
$ws.Range("D32").Value = "'0.09239"
$ws.Range("E32").Value = '  +2.48%  '
$ws.Range("D33").Value = "'0.05105"
$ws.Range("E33").Value = '  -0.39%  '
$ws.Range("D34").Value = "'0.7527"
$ws.Range("E34").Value = '  +2.57%  '
$ws.Range("D35").Value = "'2.995"
$ws.Range("E35").Value = '  -2.88%  '
$ws.Range("D36").Value = "'1.149"
$ws.Range("E36").Value = '  -0.82%  '
$ws.Range("D37").Value = "'3.289"
$ws.Range("E37").Value = '  +7.48%  '
$ws.Range("D38").Value = "'2.530"
$ws.Range("E38").Value = '  +0.65%  '
$ws.Range("D39").Value = "'0.02002"
$ws.Range("E39").Value = '  -2.25%  '
$ws.Range("D40").Value = "'0.5545"
$ws.Range("E40").Value = '  +4.35%  '
$ws.Range("D41").Value = "'1.070"
$ws.Range("E41").Value = '  -0.79%  '
$ws.Range("D42").Value = "'118.39"
$ws.Range("E42").Value = '  +1.57%  '
$ws.Range("D43").Value = "'6.542"
$ws.Range("E43").Value = '  -0.27%  '
$ws.Range("D44").Value = "'8.506"
$ws.Range("E44").Value = '  +2.18%  '
$ws.Range("D45").Value = "'0.1469"
$ws.Range("E45").Value = '  +0.07%  '
$ws.Range("D46").Value = "'0.4692"
$ws.Range("E46").Value = '  +1.50%  '
$ws.Range("D47").Value = "'0.9991"
$ws.Range("E47").Value = '  +0.04%  '
$ws.Range("D48").Value = "'10.03"
$ws.Range("E48").Value = '  -0.15%  '
$ws.Range("D49").Value = "'1.565"
$ws.Range("E49").Value = '  -0.15%  '
$ws.Range("D50").Value = "'36.78"
$ws.Range("E50").Value = '  +0.94%  '
$ws.Range("D51").Value = "'63.16"
$ws.Range("E51").Value = '  -1.92%  '
